# Apply ranking reorder fix to "max-arrecad" and "tx-sucesso" sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet "max-arrecad": rows 2-5 in column A get rotated ---
$wsMax = $wb.Worksheets.Item("max-arrecad")
$wsMax.Range("A2").Value = "fantasia"
$wsMax.Range("A3").Value = "questoes_genero"
$wsMax.Range("A4").Value = "fiq"
$wsMax.Range("A5").Value = "ficcao_cientifica"

# --- Sheet "tx-sucesso": rows 4-5 and rows 8-9 in column A get swapped ---
$wsTx = $wb.Worksheets.Item("tx-sucesso")
$wsTx.Range("A4").Value = "saloes_humor"
$wsTx.Range("A5").Value = "angelo_agostini"
$wsTx.Range("A8").Value = "questoes_genero"
$wsTx.Range("A9").Value = "erotismo"
